$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "92.733.69"
$ws.Range("E2").Value = "  -2.04%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.406.13"
$ws.Range("E3").Value = "  -0.75%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.08"
$ws.Range("E5").Value = "  -4.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "617.65"
$ws.Range("E6").Value = "  -4.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.36"
$ws.Range("E7").Value = "  -5.95%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.389"
$ws.Range("E8").Value = "  -4.46%  "

$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.959"
$ws.Range("E10").Value = "  -2.62%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.405.41"
$ws.Range("E11").Value = "  -0.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.88"
$ws.Range("E12").Value = "  +1.46%  "

$ws.Range("E13").Value = "  -1.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.21"
$ws.Range("E14").Value = "  -0.76%  "

$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.618.14"
$ws.Range("E15").Value = "  -1.92%  "

$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.043.24"
$ws.Range("E16").Value = "  -0.65%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000244"
$ws.Range("E17").Value = "  -3.64%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.11"
$ws.Range("E18").Value = "  -3.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.406.42"
$ws.Range("E19").Value = "  -0.61%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.78"
$ws.Range("E20").Value = "  +0.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.46"
$ws.Range("E21").Value = "  -2.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "493.62"
$ws.Range("E22").Value = "  -1.60%  "

$ws.Range("E23").Value = "  +1.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.430"
$ws.Range("E24").Value = "  -13.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.50"
$ws.Range("E25").Value = "  -1.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000183"
$ws.Range("E26").Value = "  -5.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "94.46"
$ws.Range("E27").Value = "  -0.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.590.13"
$ws.Range("E28").Value = "  -0.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.80"
$ws.Range("E29").Value = "  -1.95%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.28"
$ws.Range("E30").Value = "  -4.65%  "

$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("E32").Value = "  -2.48%  "

$ws.Range("E33").Value = "  -3.79%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.01"
$ws.Range("E34").Value = "  +0.41%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.171"
$ws.Range("E35").Value = "  -4.88%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "29.65"
$ws.Range("E36").Value = "  -0.75%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.537"
$ws.Range("E37").Value = "  -2.88%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "544.48"
$ws.Range("E38").Value = "  -4.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.40"
$ws.Range("E39").Value = "  -4.25%  "

$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.38"
$ws.Range("E41").Value = "  -5.66%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.913"
$ws.Range("E42").Value = "  +0.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.148"
$ws.Range("E43").Value = "  -1.73%  "

$ws.Range("E44").Value = "  -1.79%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.70"
$ws.Range("E45").Value = "  -2.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.66"
$ws.Range("E46").Value = "  -0.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.45"
$ws.Range("E47").Value = "  -4.66%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0402"
$ws.Range("E48").Value = "  -2.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.47"
$ws.Range("E49").Value = "  -3.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.09"
$ws.Range("E50").Value = "  -5.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.94"
$ws.Range("E51").Value = "  -2.10%  "
